$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting existing rows 62-114 down to 63-115.
$ws.Rows(62).Insert()

# Populate the newly inserted row 62 with the new data entry.
$ws.Cells.Item(62, 1).Value = 10
$ws.Cells.Item(62, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(62, 3).Value = "La Araucanía"
$ws.Cells.Item(62, 4).Value = 45264
$ws.Cells.Item(62, 5).Value = 9
$ws.Cells.Item(62, 6).Value = 100112042
$ws.Cells.Item(62, 7).Value = "Locoto"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 80
$ws.Cells.Item(62, 11).Value = 3300
$ws.Cells.Item(62, 12).Value = 3300
$ws.Cells.Item(62, 13).Value = 3300
$ws.Cells.Item(62, 14).Value = "$/kilo"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 3300
$ws.Cells.Item(62, 17).Value = 1
$ws.Cells.Item(62, 18).Value = "Hortaliza"
